$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Accesorio 2".."Accesorio 7" columns (M:R) entirely - both
# header (row 1) and the now-obsolete accessory breakdown values (row 2).
$ws.Range("M1:R2").EntireColumn.Delete()

# Refresh the single data row (row 2) with the new ticket's data.
$ws.Range("A2").Value = "12/31/2025"
$ws.Range("B2").Value = 4358589
$ws.Range("E2").Value = "GPS"
$ws.Range("F2").Value = "AV. LO ESPEJO 1300, MAIPU"
$ws.Range("G2").Value = "MAIPU"
$ws.Range("H2").Value = "Región Metropolitana de Santiago."
$ws.Range("J2").Value = "GENERADOR_SANTIAGO_7"
$ws.Range("K2").Value = "LUREYE"
$ws.Range("I2").Value = "Juan Perez"

# Match Excel's recalculated "best fit" column widths for the surviving
# columns (A:L) now that their contents changed.
$ws.Range("A1").EntireColumn.ColumnWidth = 9.096354166666666
$ws.Range("B1").EntireColumn.ColumnWidth = 7.029947916666667
$ws.Range("C1").EntireColumn.ColumnWidth = 7.498697916666667
$ws.Range("D1").EntireColumn.ColumnWidth = 10.233072916666666
$ws.Range("E1").EntireColumn.ColumnWidth = 8.565104166666666
$ws.Range("F1").EntireColumn.ColumnWidth = 22.498697916666668
$ws.Range("G1").EntireColumn.ColumnWidth = 6.764322916666667
$ws.Range("H1").EntireColumn.ColumnWidth = 28.299479166666668
$ws.Range("I1").EntireColumn.ColumnWidth = 13.299479166666666
$ws.Range("J1").EntireColumn.ColumnWidth = 21.897135416666668
$ws.Range("K1").EntireColumn.ColumnWidth = 5.963541666666667
$ws.Range("L1").EntireColumn.ColumnWidth = 9.233072916666666

# Move the active selection, matching the author's final cursor position.
$ws.Range("E7").Select() | Out-Null
